# Attendance record update — "day 18 on hold"
#
# 1. Row 22 (day 18, serial 21): mark "is present" = "no"
# 2. Append row 23: serial 22, day count 19, date 18-02-2026
# 3. Append row 24: serial 23, day count 20, date 19-02-2026

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that LOOKS like a number but must be stored as TEXT,
# matching the rest of the sheet (every Serial No. / Day Count cell is text).
# A direct `.Value = "22"` assignment lets Excel's type-inference store it as
# a real number, so instead we compute it as a text formula result and then
# convert that formula to a static value in place (Copy + PasteSpecial
# values-only), which keeps the cell's existing General style/format.
function Set-TextNumber {
    param($range, [string]$text)

    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $range.NumberFormat = "General"
}

# --- Day 18 put on hold: mark attendance as "no" ---
$ws.Range("D22").Value = "no"

# --- New row 23 ---
Set-TextNumber $ws.Range("A23") "22"
Set-TextNumber $ws.Range("B23") "19"
$ws.Range("C23").Value = "18-02-2026"

# --- New row 24 ---
Set-TextNumber $ws.Range("A24") "23"
Set-TextNumber $ws.Range("B24") "20"
$ws.Range("C24").Value = "19-02-2026"
